$wb = $excel.ActiveWorkbook

# Sheet "展览" (Sheet1): update 想去人数 (F column) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 3263
$ws1.Range("F3").Value = 8
$ws1.Range("F4").Value = 55
$ws1.Range("F5").Value = 1194
$ws1.Range("F6").Value = 310

# Sheet "全部类型" (Sheet4): update 想去人数 (F column) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 3263
$ws4.Range("F3").Value = 8
$ws4.Range("F4").Value = 55
$ws4.Range("F5").Value = 1194
$ws4.Range("F7").Value = 310
